$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 2
    4  = -1
    6  = -2
    7  = -3
    8  = -1
    9  = -3
    10 = -4
    11 = -4
    12 = -2
    13 = 7
    14 = 8
    15 = 2
    16 = -2
    17 = 4
    19 = -1
    20 = -2
    21 = -1
    22 = -2
    24 = -5
    25 = -1
    26 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
